# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" (Sheet1) and "全部类型" (Sheet4).
# Both sheets carry identical data and both need identical updates.

$wb = $excel.ActiveWorkbook

# Row -> new F value
$updates = @{
    2  = 1072
    3  = 772
    4  = 265
    5  = 34
    8  = 1808
    9  = 6766
    10 = 491
    11 = 382
    12 = 317
    13 = 111
    14 = 383
    15 = 142
    16 = 6972
    18 = 1306
    22 = 118
    23 = 283
    24 = 116
    27 = 103
    29 = 396
    32 = 83
    34 = 64
    36 = 65
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
